$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in week's notes, in the order that reproduces the shared-string table ---
# 1st new string: "Keras Saver fixen" (Monday / A2, also reused in C2)
$ws.Cells.Item(2, 1).Value = "Keras Saver fixen"
$ws.Cells.Item(2, 3).Value = "Keras Saver fixen"

# 2nd new string: "Aan document werken" (re-used across E2, A3:D3)
$ws.Cells.Item(2, 5).Value = "Aan document werken"
$ws.Cells.Item(3, 1).Value = "Aan document werken"
$ws.Cells.Item(3, 2).Value = "Aan document werken"
$ws.Cells.Item(3, 3).Value = "Aan document werken"
$ws.Cells.Item(3, 4).Value = "Aan document werken"

# 3rd new string: "Te veel memory van gafische kaart fixen" (B4)
$ws.Cells.Item(4, 2).Value = "Te veel memory van gafische kaart fixen"

# 4th new string: "Keras/Tensorflow Saver fixen" (B2)
$ws.Cells.Item(2, 2).Value = "Keras/Tensorflow Saver fixen"

# 5th new string: "Keras Saver fixen (het werkt!)" (D2)
$ws.Cells.Item(2, 4).Value = "Keras Saver fixen (het werkt!)"

# --- Re-style cells that previously used the now-unused "Calibri Light" xfs ---
# Copy formatting from A2 (style index 2, vertical-center/border) onto B2, C2, D2
$ws.Cells.Item(2, 1).Copy()
$ws.Range($ws.Cells.Item(2, 2), $ws.Cells.Item(2, 4)).PasteSpecial(-4122)

# Copy formatting from A3 (style index 3, border-only) onto E2, C3, D3
$ws.Cells.Item(3, 1).Copy()
$ws.Cells.Item(2, 5).PasteSpecial(-4122)
$ws.Range($ws.Cells.Item(3, 3), $ws.Cells.Item(3, 4)).PasteSpecial(-4122)

$excel.CutCopyMode = $false

# E3 is fully cleared (cell removed, not just blanked)
$ws.Cells.Item(3, 5).Clear()

# --- Update selection to match the saved view state ---
$ws.Range("D5").Select()
